$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.618.71'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '3.221.35'
$ws.Range("E3").Value = '  +0.86%  '

$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("E5").Value = '  +2.13%  '

$c = $ws.Range("D6")
$c.Value = "'158.66"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.92%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.219.85'
$ws.Range("E8").Value = '  +0.91%  '

$c = $ws.Range("D9")
$c.Value = "'0.551"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.56%  '

$ws.Range("E10").Value = '  +0.48%  '

$c = $ws.Range("D11")
$c.Value = "'5.71"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -4.96%  '

$ws.Range("E12").Value = '  -2.76%  '

$ws.Range("E13").Value = '  +1.03%  '

$c = $ws.Range("D14")
$c.Value = "'38.80"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.63%  '

$ws.Range("D15").Value = '3.751.44'
$ws.Range("E15").Value = '  +0.79%  '

$ws.Range("D16").Value = '66.656.36'
$ws.Range("E16").Value = '  +0.38%  '

$c = $ws.Range("D17")
$c.Value = "'7.36"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.22%  '

$ws.Range("D18").Value = '3.224.94'
$ws.Range("E18").Value = '  +0.91%  '

$ws.Range("E19").Value = '  +1.12%  '

$c = $ws.Range("D20")
$c.Value = "'507.34"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.46%  '

$ws.Range("E21").Value = '  -1.44%  '

$ws.Range("E22").Value = '  -0.93%  '

$c = $ws.Range("D23")
$c.Value = "'8.00"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.03%  '

$c = $ws.Range("D24")
$c.Value = "'14.61"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.73%  '

$ws.Range("E25").Value = '  -1.01%  '

$ws.Range("E26").Value = '  +0.16%  '

$ws.Range("E27").Value = '  -0.06%  '

$c = $ws.Range("D28")
$c.Value = "'9.12"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.76%  '

$c = $ws.Range("D29")
$c.Value = "'2.36"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.35%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D30")
$c.Value = "'0.121"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +34.42%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D31")
$c.Value = "'7.02"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.39%  '

$c = $ws.Range("D32")
$c.Value = "'2.94"
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.97%  '

$c = $ws.Range("D33")
$c.Value = "'28.12"
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.60%  '

$ws.Range("E34").Value = '  +0.11%  '

$ws.Range("E35").Value = '  -3.78%  '

$ws.Range("E36").Value = '  -1.08%  '

$c = $ws.Range("D37")
$c.Value = "'55.48"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.10%  '

$c = $ws.Range("D38")
$c.Value = "'500.91"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.62%  '

$ws.Range("D39").Value = '0.0₃0773'
$ws.Range("E39").Value = '  +14.96%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D40")
$c.Value = "'0.132"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +6.89%  '

$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D41")
$c.Value = "'3.08"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +7.09%  '

$c = $ws.Range("D42")
$c.Value = "'0.0420"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.09%  '

$c = $ws.Range("D43")
$c.Value = "'8.70"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("E44").Value = '  -1.74%  '

$c = $ws.Range("D45")
$c.Value = "'2.46"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.29%  '

$ws.Range("D46").Value = '2.896.39'
$ws.Range("E46").Value = '  -0.73%  '

$c = $ws.Range("D47")
$c.Value = "'28.17"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.85%  '

$ws.Range("E48").Value = '  +3.07%  '

$ws.Range("E50").Value = '  -1.10%  '

$c = $ws.Range("D51")
$c.Value = "'122.19"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.72%  '
